$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84
$ws.Range("A84").Value = "2024-01-13 22:48:13"
$ws.Range("B84").Value = 18
$ws.Range("C84").Value = 18
$ws.Range("D84").Value = 5
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 3
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0.001
$ws.Range("J84").Value = 0.05
$ws.Range("K84").Value = 0.003
$ws.Range("L84").Value = 100
$ws.Range("M84").Value = 500
$ws.Range("N84").Value = 10
$ws.Range("O84").Value = 5
$ws.Range("P84").Value = ""
$ws.Range("Q84").Value = "Data/bombay1.xlsx"

# Row 85
$ws.Range("A85").Value = "2024-01-13 23:28:28"
$ws.Range("B85").Value = 0
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0.001
$ws.Range("J85").Value = 0.05
$ws.Range("K85").Value = 0.003
$ws.Range("L85").Value = 100
$ws.Range("M85").Value = 500
$ws.Range("N85").Value = 10
$ws.Range("O85").Value = 5
$ws.Range("P85").Value = 0
$ws.Range("Q85").Value = "Data/lighting.xlsx"
